$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1.732704007046913
$ws.Range("C1").Value = 0.1945306715051764
$ws.Range("D1").Value = -0.672838158291254
$ws.Range("E1").Value = 0.7034277224914169
$ws.Range("F1").Value = 1.570796292848413
$ws.Range("G1").Value = 0.161907685808031

$ws.Range("A2").Value = 0.2074001699739488
$ws.Range("B2").Value = 1.734449897926843
$ws.Range("C2").Value = 0.1944889771920187
$ws.Range("D2").Value = -0.6729541981860235
$ws.Range("E2").Value = 0.7033533778325018
$ws.Range("F2").Value = 1.570796292708984
$ws.Range("G2").Value = 0.1636535767152294

$ws.Range("A3").Value = 0.4148003399478976
$ws.Range("B3").Value = 1.745338072090327
$ws.Range("C3").Value = 0.1942289523826182
$ws.Range("D3").Value = -0.6736778760508239
$ws.Range("E3").Value = 0.7028897305314603
$ws.Range("F3").Value = 1.570796291839441
$ws.Range("G3").Value = 0.1745417510487718

$ws.Range("A4").Value = 0.6222005099218464
$ws.Range("B4").Value = 1.771088086074367
$ws.Range("C4").Value = 0.1936140060768625
$ws.Range("D4").Value = -0.6753893397782359
$ws.Range("E4").Value = 0.7017932267185878
$ws.Range("F4").Value = 1.570796289783013
$ws.Range("G4").Value = 0.2002917654349918

$ws.Range("A5").Value = 0.8296006798957951
$ws.Range("B5").Value = 1.814176913182388
$ws.Range("C5").Value = 0.1925849846920256
$ws.Range("D5").Value = -0.6782532203453618
$ws.Range("E5").Value = 0.6999583903085791
$ws.Range("F5").Value = 1.570796286341886
$ws.Range("G5").Value = 0.2433805932160021

$ws.Range("A6").Value = 1.037000849869744
$ws.Range("B6").Value = 1.87437937402317
$ws.Range("C6").Value = 0.1911472658265371
$ws.Range("D6").Value = -0.6822545512997377
$ws.Range("E6").Value = 0.6973948100364618
$ws.Range("F6").Value = 1.570796281534042
$ws.Range("G6").Value = 0.3035830549970654

$ws.Range("A7").Value = 1.244401019843693
$ws.Range("B7").Value = 1.949308567049777
$ws.Range("C7").Value = 0.1893578520237512
$ws.Range("D7").Value = -0.6872346882452469
$ws.Range("E7").Value = 0.6942041264935301
$ws.Range("F7").Value = 1.570796275550103
$ws.Range("G7").Value = 0.3785122491939638

$ws.Range("A8").Value = 1.451801189817642
$ws.Range("B8").Value = 2.03495629909848
$ws.Range("C8").Value = 0.1873124645357151
$ws.Range("D8").Value = -0.6929272283280338
$ws.Range("E8").Value = 0.6905570191632778
$ws.Range("F8").Value = 1.570796268710169
$ws.Range("G8").Value = 0.4641599825803687

$ws.Range("A9").Value = 1.65920135979159
$ws.Range("B9").Value = 2.126233515927697
$ws.Range("C9").Value = 0.1851326370869383
$ws.Range("D9").Value = -0.6989939297224161
$ws.Range("E9").Value = 0.686670193457332
$ws.Range("F9").Value = 1.570796261420657
$ws.Range("G9").Value = 0.5554372008352122

$ws.Range("A10").Value = 1.866601529765539
$ws.Range("B10").Value = 2.217510732756915
$ws.Range("C10").Value = 0.1829528096381616
$ws.Range("D10").Value = -0.7050606311167984
$ws.Range("E10").Value = 0.6827833677513863
$ws.Range("F10").Value = 1.570796254131145
$ws.Range("G10").Value = 0.6467144190900558

$ws.Range("A11").Value = 2.074001699739488
$ws.Range("B11").Value = 2.303158464805618
$ws.Range("C11").Value = 0.1809074221501255
$ws.Range("D11").Value = -0.7107531711995853
$ws.Range("E11").Value = 0.679136260421134
$ws.Range("F11").Value = 1.57079624729121
$ws.Range("G11").Value = 0.732362152476461

$ws.Range("A12").Value = 2.281401869713437
$ws.Range("B12").Value = 2.378087657832224
$ws.Range("C12").Value = 0.1791180083473395
$ws.Range("D12").Value = -0.7157333081450945
$ws.Range("E12").Value = 0.6759455768782023
$ws.Range("F12").Value = 1.570796241307272
$ws.Range("G12").Value = 0.8072913466733589

$ws.Range("A13").Value = 2.488802039687386
$ws.Range("B13").Value = 2.438290118673007
$ws.Range("C13").Value = 0.1776802894818511
$ws.Range("D13").Value = -0.7197346390994703
$ws.Range("E13").Value = 0.673381996606085
$ws.Range("F13").Value = 1.570796236499428
$ws.Range("G13").Value = 0.8674938084544226

$ws.Range("A14").Value = 2.696202209661334
$ws.Range("B14").Value = 2.481378945781028
$ws.Range("C14").Value = 0.1766512680970142
$ws.Range("D14").Value = -0.7225985196665963
$ws.Range("E14").Value = 0.6715471601960763
$ws.Range("F14").Value = 1.570796233058301
$ws.Range("G14").Value = 0.9105826362354328

$ws.Range("A15").Value = 2.903602379635283
$ws.Range("B15").Value = 2.507128959765068
$ws.Range("C15").Value = 0.1760363217912584
$ws.Range("D15").Value = -0.7243099833940083
$ws.Range("E15").Value = 0.6704506563832038
$ws.Range("F15").Value = 1.570796231001873
$ws.Range("G15").Value = 0.9363326506216527

$ws.Range("A16").Value = 3.111002549609232
$ws.Range("B16").Value = 2.518017133928552
$ws.Range("C16").Value = 0.1757762969818579
$ws.Range("D16").Value = -0.7250336612588086
$ws.Range("E16").Value = 0.6699870090821622
$ws.Range("F16").Value = 1.57079623013233
$ws.Range("G16").Value = 0.947220824955195

$ws.Range("A17").Value = 3.31840271958318
$ws.Range("B17").Value = 2.519763024808483
$ws.Range("C17").Value = 0.1757346026687003
$ws.Range("D17").Value = -0.7251497011535781
$ws.Range("E17").Value = 0.6699126644232472
$ws.Range("F17").Value = 1.570796229992901
$ws.Range("G17").Value = 0.9489667158623929

